$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

function Set-TextCell($r, $c, $text) {
    # Leading apostrophe forces a literal-text entry (prevents Excel from
    # auto-converting things like "55%" into a formatted percentage number).
    # Resetting the style back to Normal afterwards drops the transient
    # "quote prefix" style flag that the apostrophe entry leaves behind, so
    # the cell ends up with the default style (matching the rest of the
    # sheet) while still holding the literal text value.
    $ws.Cells.Item($r, $c).Value = "'" + $text
    $ws.Cells.Item($r, $c).Style = "Normal"
}

function Set-EmptyTextCell($r, $c) {
    Set-TextCell $r $c ""
}

function Set-NumberCell($r, $c, $num) {
    $ws.Cells.Item($r, $c).Value = $num
}

Set-TextCell $row 1 "2024-10-20 01:59:54"
Set-EmptyTextCell $row 2
Set-NumberCell $row 3 12
Set-NumberCell $row 4 1
Set-NumberCell $row 5 3
Set-NumberCell $row 6 8
Set-NumberCell $row 7 0
Set-NumberCell $row 8 0
Set-EmptyTextCell $row 9
Set-EmptyTextCell $row 10
Set-EmptyTextCell $row 11
Set-EmptyTextCell $row 12
Set-EmptyTextCell $row 13
Set-NumberCell $row 14 10
Set-NumberCell $row 15 10
Set-NumberCell $row 16 2
Set-EmptyTextCell $row 17
Set-NumberCell $row 18 5
Set-EmptyTextCell $row 19
Set-NumberCell $row 20 20
Set-TextCell $row 21 "55%"
Set-TextCell $row 22 "C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx"
Set-EmptyTextCell $row 23
Set-TextCell $row 24 "No es Simulación"
Set-NumberCell $row 25 22
